$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 193, shifting existing rows 193:298 down to 194:299
$ws.Rows("193:193").Insert()

# Populate the newly inserted row 193 with the new weekly data point
$ws.Range("A193").Value = 4
$ws.Range("B193").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C193").Value = "Los Lagos"
$ws.Range("D193").Value = 44719
$ws.Range("D193").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E193").Value = 10
$ws.Range("F193").Value = 100112040
$ws.Range("G193").Value = "Cilantro"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 180
$ws.Range("K193").Value = 11000
$ws.Range("L193").Value = 11000
$ws.Range("M193").Value = 11000
$ws.Range("N193").Value = "`$/caja 36 atados"
$ws.Range("O193").Value = "Región Metropolitana"
$ws.Range("P193").Value = 306
$ws.Range("Q193").Value = 36
$ws.Range("R193").Value = "Hortaliza"
